$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$headers = @("code","name","day_seq","is_global_working","lang_code","is_active","cr_by","cr_dtimes","upd_by","upd_dtimes","is_deleted","del_dtimes")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---- Number format for the timestamp columns (H = cr_dtimes) ----
$ws.Range("H2:H8").NumberFormat = "mm:ss.0"

# ---- Data rows ----
# columns: A=code  B=name  C=day_seq  D=is_global_working  E=lang_code  F=is_active
#          G=cr_by H=cr_dtimes I=upd_by J=upd_dtimes K=is_deleted L=del_dtimes
$data = @(
    @(101, "DIM", 1, $false, "fra", $true),
    @(102, "LUN", 2, $true,  "fra", $true),
    @(103, "MAR", 3, $true,  "fra", $true),
    @(104, "MER", 4, $true,  "fra", $true),
    @(105, "JEU", 5, $true,  "fra", $true),
    @(106, "VEN", 6, $true,  "fra", $true),
    @(107, "SAM", 7, $false, "fra", $true)
)

$cr_dtimes = 45079.57763521991

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = "superadmin"
    $ws.Cells.Item($r, 8).Value = $cr_dtimes
    $ws.Cells.Item($r, 9).Value = "NULL"
    $ws.Cells.Item($r, 10).Value = "NULL"
    $ws.Cells.Item($r, 11).Value = $false
    $ws.Cells.Item($r, 12).Value = "NULL"
}

# ---- Selection, to mirror the saved cursor position in the target file ----
[void]$ws.Range("D15").Select()
